$wb = $excel.ActiveWorkbook

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/17e2bf67959ea0264658d4a9c09fe1849e060cb7/e2e/a.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/17e2bf67959ea0264658d4a9c09fe1849e060cb7/e2e/b.md"

# ---------------------------------------------------------------
# Overview sheet: status text update + widened zh-cn/de-de columns
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Columns(5).ColumnWidth = 29.09
$wsOverview.Columns(6).ColumnWidth = 29.09

# ---------------------------------------------------------------
# zh-cn sheet: handback file/date populated, column widths widened
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns(3).ColumnWidth = 29.09
$wsZh.Columns(10).ColumnWidth = 39.17

$wsZh.Range("I2").Value = "a.md"
$wsZh.Range("I2").Style = "Hyperlink"
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-31 14:46:48"

$wsZh.Range("I3").Value = "a.md"
$wsZh.Range("I3").Style = "Hyperlink"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-31 14:46:48"

$wsZh.Range("A3").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlA, "", "", "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, "", "", "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlB, "", "", "b.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlA, "", "", "a.md") | Out-Null

# ---------------------------------------------------------------
# de-de sheet: handback file/date populated, column widths widened
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns(3).ColumnWidth = 29.09
$wsDe.Columns(10).ColumnWidth = 39.17

$wsDe.Range("I2").Value = "a.md"
$wsDe.Range("I2").Style = "Hyperlink"
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-31 14:46:55"

$wsDe.Range("I3").Value = "a.md"
$wsDe.Range("I3").Style = "Hyperlink"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-31 14:46:55"

$wsDe.Range("A3").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlA, "", "", "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, "", "", "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlB, "", "", "b.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlA, "", "", "a.md") | Out-Null
